# Update: Fix new logic for import data
# Insert a new "Merk" column before the existing "Jenis" column (E),
# shifting the later columns right, and append "User"/"Dept" columns
# at the end of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new column E ("Merk") before old column E ("Jenis") ---
$ws.Columns.Item(5).Insert()
$ws.Range("E1").Value = "Merk"

# --- Append the two new trailing columns ---
$ws.Range("L1").Value = "User"
$ws.Range("M1").Value = "Dept"

# The header row style (bold white text on blue fill) needs to be copied
# onto the two newly appended header cells - a plain .Value assignment
# does not carry any formatting with it.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
[void]$ws.Application.CutCopyMode

# --- Column widths (best effort via the character-based ColumnWidth API) ---
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 14.333333333333334
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668
$ws.Columns.Item(7).ColumnWidth = 16
$ws.Columns.Item(8).ColumnWidth = 19.666666666666668
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(11).ColumnWidth = 21
$ws.Columns.Item(12).ColumnWidth = 21
$ws.Columns.Item(13).ColumnWidth = 17

# --- Selection moved by the author while reviewing the new columns ---
[void]$ws.Range("I17").Select()
